$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.840.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.257.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.40%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.450"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.597.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.40%  "
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.263.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.917.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -4.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +22.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("E34").Value = "  +7.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  -3.26%  "
$ws.Range("E37").Value = "  +4.41%  "
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0261"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.43%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0967"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("E45").Value = "  -5.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "97.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("E49").Value = "  +6.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.439.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.35%  "
